$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change) per latest scrape.
$ws.Cells.Item(2, 4).Value = '25.602.73'
$ws.Cells.Item(2, 5).Value = '  +2.26%  '
$ws.Cells.Item(3, 4).Value = '1.665.55'
$ws.Cells.Item(3, 5).Value = '  +0.81%  '
$ws.Cells.Item(4, 4).Value = '''0.9992'
$ws.Cells.Item(4, 5).Value = '  -0.20%  '
$ws.Cells.Item(5, 4).Value = '''236.70'
$ws.Cells.Item(5, 5).Value = '  +0.04%  '
$ws.Cells.Item(6, 5).Value = '  -0.09%  '
$ws.Cells.Item(7, 4).Value = '''0.4806'
$ws.Cells.Item(7, 5).Value = '  +0.58%  '
$ws.Cells.Item(8, 4).Value = '''0.2630'
$ws.Cells.Item(8, 5).Value = '  +0.46%  '
$ws.Cells.Item(9, 4).Value = '''0.06158'
$ws.Cells.Item(9, 5).Value = '  +2.80%  '
$ws.Cells.Item(10, 4).Value = '''0.07087'
$ws.Cells.Item(10, 5).Value = '  -0.02%  '
$ws.Cells.Item(11, 4).Value = '1.663.92'
$ws.Cells.Item(11, 5).Value = '  +0.60%  '
$ws.Cells.Item(12, 5).Value = '  +2.87%  '
$ws.Cells.Item(13, 4).Value = '''0.5923'
$ws.Cells.Item(13, 5).Value = '  -4.11%  '
$ws.Cells.Item(14, 4).Value = '''4.393'
$ws.Cells.Item(14, 5).Value = '  -4.02%  '
$ws.Cells.Item(15, 4).Value = '''74.50'
$ws.Cells.Item(15, 5).Value = '  +2.15%  '
$ws.Cells.Item(16, 5).Value = '  -0.04%  '
$ws.Cells.Item(17, 4).Value = '''0.9997'
$ws.Cells.Item(17, 5).Value = '  -0.21%  '
$ws.Cells.Item(18, 4).Value = '25.581.56'
$ws.Cells.Item(18, 5).Value = '  +2.17%  '
$ws.Cells.Item(19, 4).Value = '''0.000006763'
$ws.Cells.Item(19, 5).Value = '  +2.91%  '
$ws.Cells.Item(20, 5).Value = '  +0.64%  '
$ws.Cells.Item(21, 4).Value = '1.879.75'
$ws.Cells.Item(21, 5).Value = '  +0.75%  '
$ws.Cells.Item(22, 5).Value = '  +0.84%  '
$ws.Cells.Item(23, 4).Value = '''8.681'
$ws.Cells.Item(23, 5).Value = '  +2.77%  '
$ws.Cells.Item(24, 4).Value = '''5.313'
$ws.Cells.Item(24, 5).Value = '  +1.43%  '
$ws.Cells.Item(25, 4).Value = '''134.88'
$ws.Cells.Item(26, 4).Value = '''15.06'
$ws.Cells.Item(26, 5).Value = '  +2.39%  '
$ws.Cells.Item(27, 4).Value = '''1.405'
$ws.Cells.Item(27, 5).Value = '  +0.99%  '
$ws.Cells.Item(28, 4).Value = '''105.02'
$ws.Cells.Item(28, 5).Value = '  +3.56%  '
$ws.Cells.Item(29, 4).Value = '''1.692'
$ws.Cells.Item(29, 5).Value = '  +0.06%  '
$ws.Cells.Item(30, 4).Value = '''3.953'
$ws.Cells.Item(30, 5).Value = '  +4.73%  '
$ws.Cells.Item(31, 4).Value = '''3.671'
$ws.Cells.Item(31, 5).Value = '  +4.24%  '
$ws.Cells.Item(32, 4).Value = '''0.07659'
$ws.Cells.Item(32, 5).Value = '  -3.13%  '
$ws.Cells.Item(33, 4).Value = '''0.9996'
$ws.Cells.Item(34, 4).Value = '''0.04324'
$ws.Cells.Item(34, 5).Value = '  -5.30%  '
$ws.Cells.Item(35, 4).Value = '''2.615'
$ws.Cells.Item(35, 5).Value = '  -0.43%  '
$ws.Cells.Item(36, 4).Value = '''0.6135'
$ws.Cells.Item(36, 5).Value = '  +6.27%  '
$ws.Cells.Item(37, 4).Value = '''0.9509'
$ws.Cells.Item(37, 5).Value = '  +1.44%  '
$ws.Cells.Item(38, 4).Value = '''2.609'
$ws.Cells.Item(38, 5).Value = '  -0.49%  '
$ws.Cells.Item(39, 4).Value = '''0.8595'
$ws.Cells.Item(39, 5).Value = '  +1.91%  '
$ws.Cells.Item(40, 5).Value = '  -0.08%  '
$ws.Cells.Item(41, 5).Value = '  -1.42%  '
$ws.Cells.Item(42, 5).Value = '  +2.90%  '
$ws.Cells.Item(43, 4).Value = '''98.05'
$ws.Cells.Item(43, 5).Value = '  -0.49%  '
$ws.Cells.Item(44, 4).Value = '''0.3771'
$ws.Cells.Item(44, 5).Value = '  +2.16%  '
$ws.Cells.Item(45, 5).Value = '  -2.44%  '
$ws.Cells.Item(46, 4).Value = '''0.1121'
$ws.Cells.Item(46, 5).Value = '  +0.89%  '
$ws.Cells.Item(47, 4).Value = '''6.231'
$ws.Cells.Item(47, 5).Value = '  +3.04%  '
$ws.Cells.Item(48, 4).Value = '''0.05263'
$ws.Cells.Item(49, 4).Value = '''29.51'
$ws.Cells.Item(49, 5).Value = '  +0.68%  '
$ws.Cells.Item(50, 4).Value = '''7.375'
$ws.Cells.Item(50, 5).Value = '  +1.41%  '
$ws.Cells.Item(51, 4).Value = '''1.001'
$ws.Cells.Item(51, 5).Value = '  -0.01%  '
